$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Simple D/E value updates ---
Set-TextCell $ws.Range("D2") "256.71"
Set-TextCell $ws.Range("E2") "-0.03%"
Set-TextCell $ws.Range("E3") "-2.89%"
Set-TextCell $ws.Range("D4") "4.624"
Set-TextCell $ws.Range("E4") "-11.45%"
Set-TextCell $ws.Range("D5") "0.05895"
Set-TextCell $ws.Range("E5") "0.10%"
Set-TextCell $ws.Range("D6") "6.643"
Set-TextCell $ws.Range("E6") "-0.90%"
Set-TextCell $ws.Range("E7") "-0.70%"
Set-TextCell $ws.Range("D8") "0.9335"
Set-TextCell $ws.Range("E8") "-4.38%"
Set-TextCell $ws.Range("D16") "0.006124"
Set-TextCell $ws.Range("E16") "1.54%"
Set-TextCell $ws.Range("D17") "3.517"
Set-TextCell $ws.Range("E17") "0.52%"
Set-TextCell $ws.Range("D18") "3.188"
Set-TextCell $ws.Range("E18") "-0.90%"
Set-TextCell $ws.Range("E19") "0.32%"
Set-TextCell $ws.Range("E21") "-1.14%"
Set-TextCell $ws.Range("D22") "3.862"
Set-TextCell $ws.Range("E22") "9.46%"
Set-TextCell $ws.Range("D23") "0.04239"
Set-TextCell $ws.Range("E23") "1.69%"
Set-TextCell $ws.Range("D24") "0.001218"
Set-TextCell $ws.Range("E24") "-0.60%"
Set-TextCell $ws.Range("D25") "0.004278"
Set-TextCell $ws.Range("E25") "-6.22%"
Set-TextCell $ws.Range("D26") "0.0001199"
Set-TextCell $ws.Range("E26") "0.09%"
Set-TextCell $ws.Range("E27") "31.58%"
Set-TextCell $ws.Range("D40") "0.03830"
Set-TextCell $ws.Range("E40") "-0.11%"
Set-TextCell $ws.Range("D43") "0.002429"
Set-TextCell $ws.Range("E43") "3.74%"
Set-TextCell $ws.Range("D44") "0.01136"
Set-TextCell $ws.Range("E44") "19.00%"
Set-TextCell $ws.Range("D45") "0.00005470"
Set-TextCell $ws.Range("E45") "1.19%"
Set-TextCell $ws.Range("E46") "0.06%"
Set-TextCell $ws.Range("D47") "0.07774"
Set-TextCell $ws.Range("E47") "-18.09%"
Set-TextCell $ws.Range("D48") "0.002279"
Set-TextCell $ws.Range("E48") "7.08%"
Set-TextCell $ws.Range("E49") "0.06%"
Set-TextCell $ws.Range("D50") "0.0001998"
Set-TextCell $ws.Range("E50") "0.06%"

# --- Row 9-15 cyclic reshuffle (coin list reorder) ---
Set-TextCell $ws.Range("B9") "WazirX"
Set-TextCell $ws.Range("C9") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell $ws.Range("D9") "0.1406"
Set-TextCell $ws.Range("E9") "-0.35%"
Set-TextCell $ws.Range("B10") "LiechtensteinCryptoassetsExchange"
Set-TextCell $ws.Range("C10") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell $ws.Range("D10") "0.03845"
Set-TextCell $ws.Range("E10") "10.14%"
Set-TextCell $ws.Range("B11") "MandalaExchangeToken"
Set-TextCell $ws.Range("C11") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell $ws.Range("D11") "0.07086"
Set-TextCell $ws.Range("E11") "-1.47%"
Set-TextCell $ws.Range("B12") "BitrueCoin"
Set-TextCell $ws.Range("C12") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell $ws.Range("D12") "0.03202"
Set-TextCell $ws.Range("E12") "0.94%"
Set-TextCell $ws.Range("B13") "BitMartToken"
Set-TextCell $ws.Range("C13") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell $ws.Range("D13") "0.09255"
Set-TextCell $ws.Range("E13") "0.38%"
Set-TextCell $ws.Range("B14") "BitForexToken"
Set-TextCell $ws.Range("C14") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell $ws.Range("D14") "0.001552"
Set-TextCell $ws.Range("E14") "0.55%"
Set-TextCell $ws.Range("B15") "One"
Set-TextCell $ws.Range("C15") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell $ws.Range("D15") "0.0006025"
Set-TextCell $ws.Range("E15") "-0.47%"

# --- Row 41-42 swap ---
Set-TextCell $ws.Range("B41") "BKEXToken"
Set-TextCell $ws.Range("C41") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell $ws.Range("D41") "0.1101"
Set-TextCell $ws.Range("E41") "-0.23%"
Set-TextCell $ws.Range("B42") "KickToken"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell $ws.Range("D42") "0.006187"
Set-TextCell $ws.Range("E42") "13.06%"
